$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the leading "Reworking Tree:" section entirely (the header
#    paragraph plus the whole numId=6 bullet list plus the blank line
#    that used to separate it from "Common plan:").
# ---------------------------------------------------------------------
$startPara = $d.Paragraphs.Item(1)
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -match "^Common plan:") {
        $endPara = $d.Paragraphs.Item($i - 1)
        break
    }
}
$killRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$killRange.Delete()

# ---------------------------------------------------------------------
# 2. Locate the "Common plan:" numId=7 sub-list and rework its first
#    four bullets, demoting three of them to a nested (ilvl=1) level
#    and inserting the new sub-bullets describing the drawing work.
# ---------------------------------------------------------------------
function Get-ParaIndexByText($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $prefix) {
            return $i
        }
    }
    return -1
}

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

function Set-ParaLevel($index, $level) {
    $d.Paragraphs.Item($index).Range.ListFormat.ListLevelNumber = $level
}

$idxFinishTree = Get-ParaIndexByText "^Finish Tree"
Set-ParaText $idxFinishTree "Implement simple drawing graphs"

$idxCover = Get-ParaIndexByText "^Cover Tree and Nodes by tests"
Set-ParaText $idxCover "Form and canvas"
Set-ParaLevel $idxCover 2

$idxDevelop = Get-ParaIndexByText "^Develop light \(or stub\) electrical-net-specific content"
Set-ParaText $idxDevelop "Block"
Set-ParaLevel $idxDevelop 2

$idxDraw = Get-ParaIndexByText "^Implement simple drawing graphs"
# the first match now is the renamed "Finish Tree" paragraph above, so
# look again starting right after it for the original 4th bullet
$idxDraw = -1
for ($i = $idxDevelop + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^Implement simple drawing graphs") {
        $idxDraw = $i
        break
    }
}
Set-ParaText $idxDraw "Simple and fool connection lines"
Set-ParaLevel $idxDraw 2

# Insert the remaining new nested bullets right after it, in order.
$newBullets = @(
    "Dragging blocks vertically",
    "Beauty connection lines",
    "Auto movement connection lines",
    "Auto aligning blocks by horizontal"
)

$anchorIndex = $idxDraw
foreach ($bulletText in $newBullets) {
    $anchorPara = $d.Paragraphs.Item($anchorIndex)
    $anchorPara.Range.InsertParagraphAfter()
    $anchorIndex = $anchorIndex + 1
    $newPara = $d.Paragraphs.Item($anchorIndex)
    $r = $newPara.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $bulletText
    $newPara.Range.ListFormat.ListLevelNumber = 2
}

# ---------------------------------------------------------------------
# 3. The nested (ilvl=1) level of that list (numId=7 / abstractNum 1)
#    was only "tentative" before it got real content; touching the
#    level definition marks it as actually used.
# ---------------------------------------------------------------------
$lt = $d.Paragraphs.Item($idxCover).Range.ListFormat.ListTemplate
$lvl = $lt.ListLevels.Item(2)
$lvl.NumberFormat = $lvl.NumberFormat
